$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.194.81'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '1.906.44'
$ws.Range("E4").Value = '  +0.19%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '307.77'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("E6").Value = '  +0.16%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5239'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +2.82%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3779'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +3.07%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.07259'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +1.08%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '21.23'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +2.79%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.8961'
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.07681'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +2.52%  '
$ws.Range("D13").Value = '1.902.66'
$ws.Range("E13").Value = '  +1.20%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '95.06'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +0.61%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '5.273'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("E16").Value = '  +0.22%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.000008591'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.98%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '14.42'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '27.262.28'
$ws.Range("E20").Value = '  +1.12%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '5.064'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("D22").Value = '2.154.20'
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("E23").Value = '  +2.47%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '6.429'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.63%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.308'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +10.75%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '145.80'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -1.49%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '1.740'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -2.01%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '18.14'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.52%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '114.73'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +1.03%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '4.966'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +5.28%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '4.798'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +2.10%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.09216'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +0.59%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.8113'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +8.25%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.05051'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -0.06%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.240'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +7.54%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '3.000'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +0.67%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '3.310'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +2.85%  '
$ws.Range("E38").Value = '  +2.51%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.5680'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +0.66%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.01984'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("E41").Value = '  +0.39%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '8.979'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +5.04%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '119.19'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +3.12%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '6.619'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.03%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.1513'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +2.09%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.4828'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +0.91%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '10.19'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("E48").Value = '  +0.17%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '1.624'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +4.17%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '37.54'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +1.51%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '63.73'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.97%  '
